$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper now also captures player "height" and "weight", inserted as two
# new columns between "fumbles" (D) and the existing "fantasy points" column
# (previously E). Move "fantasy points" out to G, then populate the freed up
# E/F columns with the new header + values.

# 1) Move "fantasy points" column (header + all 16 data rows) from E to G.
for ($r = 1; $r -le 17; $r++) {
    $src = $ws.Cells.Item($r, 5)
    $dst = $ws.Cells.Item($r, 7)
    $dst.Value = $src.Value2
}
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) New headers for the inserted columns, matching the bold/bordered header
#    style already used by the other header cells.
$ws.Cells.Item(1, 5).Value = "height"
$ws.Cells.Item(1, 6).Value = "weight"
$ws.Cells.Item(1, 4).Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate height/weight for every data row (2-17).
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.333333333333333
    $ws.Cells.Item($r, 6).Value = 254
}
